$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 'ALEX1,ALEX2,ALEX4,ALEX5,ALEX6,ALEX7,ALEX8,LAEX3,RQ1,RQ2,G17,G25,G5,G6,E,DC01,DC09,DC10,DC11,DC12,$G11,$G12,EKG1,EKG2,PST1'
$ws.Range("B7").Value = 'ALEX1,ALEX2,ALEX4,ALEX5,ALEX6,ALEX7,ALEX8,LAEX3,RQ1,RQ2,G17,G25,G5,G6,E,DC01,DC09,DC10,DC11,DC12,$G11,$G12,EKG1,EKG2,PST1'
$ws.Range("B13").Value = 'RR1,RR2,RR3,RR4,LR1,LR2,LR3,LR4,LPLT1,LPLT6,DC01,DC09,DC10,DC11,DC12,$RALT5,$RALT6,EKG1,EKG2,LFT6'
$ws.Range("B15").Value = 'G33,G34,G35,G36,G48,R1,R2,MFP5,MFP6,LFP6,LFP7,LFP8,MT6,LSF6,E,DC01,$LFP3,$LFP4,EKG1,EKG2,G42'
$ws.Range("B16").Value = 'G33,G34,G35,G36,G48,R1,R2,MFP5,MFP6,LFP6,LFP7,LFP8,MT6,LSF6,E,DC01,$LFP3,$LFP4,EKG1,EKG2,G42'
$ws.Range("B17").Value = 'MST2,PO5,PO6,PO7,PO8,PO9,PO10,R1,R2,SO7,SO8,SO9,SO10,SO11,SO12,E,DC01,$MST1,$MST2,EKG1,EKG2,MST3,PPST1,MST1,TO1'
$ws.Range("B18").Value = 'MST2,PO5,PO6,PO7,PO8,PO9,PO10,R1,R2,SO7,SO8,SO9,SO10,SO11,SO12,E,DC01,$MST1,$MST2,EKG1,EKG2,MST3,PPST1,MST1,TO1,PPST2,G25,G9,G16,G8,TO6,TO2'
$ws.Range("B19").Value = 'MST2,PO5,PO6,PO7,PO8,PO9,PO10,R1,R2,SO7,SO8,SO9,SO10,SO11,SO12,E,DC01,$MST1,$MST2,EKG1,EKG2,G8,MST1,PPST1'
$ws.Range("B23").Value = 'LIAI1,LIAI2,LIAI3,LIAI6,LIM1,LIM2,LIM3,LIM4,LIPI1,LIPI2,LIPI3,LIPI4,LIPI5,LIPI6,R1,R2,RIAS1,RIAS2,RIAS3,RIAS4,RIAI1,RIAI2,RIAI3,RIAI4,RIAI5,RIAI6,RIM1,RIM2,RIM3,RIM4,RIPI1,RIPI2,RIPI3,RIPI5,RIPI6,RIPS1,RIPS2,RIPS3,RIPS5,RIPS6,RAL8,RAM8,RAL4,RG25,E,DC01,$RG1,$RG2,EKG1,EKG2,RG3'
$ws.Range("B24").Value = 'LIAI1,LIAI2,LIAI3,LIAI6,LIM1,LIM2,LIM3,LIM4,LIPI1,LIPI2,LIPI3,LIPI4,LIPI5,LIPI6,R1,R2,RIAS1,RIAS2,RIAS3,RIAS4,RIAI1,RIAI2,RIAI3,RIAI4,RIAI5,RIAI6,RIM1,RIM2,RIM3,RIM4,RIPI1,RIPI2,RIPI3,RIPI5,RIPI6,RIPS1,RIPS2,RIPS3,RIPS5,RIPS6,RAL8,RAM8,RAL4,RG25,E,DC01,$RG1,$RG2,EKG1,EKG2,RPG9,RPG1,LG1,LG3,LIAI4-1'
$ws.Range("B25").Value = 'LIAI1,LIAI2,LIAI3,LIAI6,LIM1,LIM2,LIM3,LIM4,LIPI1,LIPI2,LIPI3,LIPI4,LIPI5,LIPI6,R1,R2,RIAS1,RIAS2,RIAS3,RIAS4,RIAI1,RIAI2,RIAI3,RIAI4,RIAI5,RIAI6,RIM1,RIM2,RIM3,RIM4,RIPI1,RIPI2,RIPI3,RIPI5,RIPI6,RIPS1,RIPS2,RIPS3,RIPS5,RIPS6,RAL8,RAM8,RAL4,RG25,E,DC01,EKG1,EKG2,$RG11,$RG12'
$ws.Range("B26").Value = 'LIAI1,LIAI2,LIAI3,LIAI6,LIM1,LIM2,LIM3,LIM4,LIPI1,LIPI2,LIPI3,LIPI4,LIPI5,LIPI6,R1,R2,RIAS1,RIAS2,RIAS3,RIAS4,RIAI1,RIAI2,RIAI3,RIAI4,RIAI5,RIAI6,RIM1,RIM2,RIM3,RIM4,RIPI1,RIPI2,RIPI3,RIPI5,RIPI6,RIPS1,RIPS2,RIPS3,RIPS5,RIPS6,RAL8,RAM8,RAL4,RG25,E,DC01,EKG1,EKG2,$RG11,$RG12,RAM1,RIPS4-0,LIPS5,LIPS1'
$ws.Range("B31").Value = 'FP1,FP2,FP3,FP4,R1,R2,RAIH3,RPIH2,RPPIH2,DC09,DC10,$G11,$G12,EKG1,EKG2,G16'
$ws.Range("B35").Value = 'G5,G6,G7,G15,G23,P1,P2,P3,P4,P5,P6,RU1,RU2,PST4,E,DC09,DC10,DC11,$G13,$G14,EKG1,EKG2,TT1'
$ws.Range("B36").Value = 'PST3,PST4,R1,R2,R3,R4,B19,LSF8,G2,G10,OF2,TT5,DC09,DC10,DC11,DC12,$G13,$G14,EKG1,EKG2,PPST3'
$ws.Range("B37").Value = 'PST3,PST4,R1,R2,R3,R4,B19,LSF8,G2,G10,OF2,TT5,DC09,DC10,DC11,DC12,$G13,$G14,EKG1,EKG2,PPST3,TT3'
$ws.Range("B38").Value = 'PST3,PST4,R1,R2,R3,R4,B19,LSF8,G2,G10,OF2,TT5,DC09,DC10,DC11,DC12,$G13,$G14,EKG1,EKG2,LIF4,LIF5'
$ws.Range("B39").Value = 'PST3,PST4,R1,R2,R3,R4,B19,LSF8,G2,G10,OF2,TT5,DC09,DC10,DC11,DC12,$G13,$G14,EKG1,EKG2,LIF6'
$ws.Range("B40").Value = 'R1,R2,R3,R4,R5,E,DC01,DC02,DC03,DC04,$AST1,$AST2,EKG1,EKG2,OF4,MST4,TT6'
$ws.Range("B41").Value = 'R1,R2,R3,R4,R5,E,DC01,DC02,DC03,DC04,$AST1,$AST2,EKG1,EKG2,MST4,TT6'

$ws.Range("B42").Select()
